# バーンダウンチャート.xlsx - update the "残作業時間" (remaining work) row (row 8)
# on sheet "スプリント１" so that each day's remaining total is derived by
# walking backwards from the grand total in column M, instead of the old
# "walk forward from day 1" shared formula.
#
#   B8 = SUM(M2:M6)            (total estimated hours)
#   C8 = B8 - SUM(C2:C6)       (remaining after day 1's work is subtracted)
#   D8:L8 = previous day's remaining - SUM(this day's hours, rows 2:6)
#
# Each cell is written individually (rather than via one multi-cell Range
# assignment) so every formula is entered exactly as a user typing them one
# by one would, instead of as one bulk fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("スプリント１")

$ws.Range("B8").Formula = "=SUM(M2:M6)"
$ws.Range("C8").Formula = "=B8-SUM(C2:C6)"
$ws.Range("D8").Formula = "=C8-SUM(D2:D6)"
$ws.Range("E8").Formula = "=D8-SUM(E2:E6)"
$ws.Range("F8").Formula = "=E8-SUM(F2:F6)"
$ws.Range("G8").Formula = "=F8-SUM(G2:G6)"
$ws.Range("H8").Formula = "=G8-SUM(H2:H6)"
$ws.Range("I8").Formula = "=H8-SUM(I2:I6)"
$ws.Range("J8").Formula = "=I8-SUM(J2:J6)"
$ws.Range("K8").Formula = "=J8-SUM(K2:K6)"
$ws.Range("L8").Formula = "=K8-SUM(L2:L6)"
